$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F13 is a stray styled cell left over from the old template; turn its style
# slot (cellXfs index 1, currently just "applyFont") into the integer number
# format that the surviving CANTIDAD DE PASAJEROS column needs, before the
# row carrying it is removed.
$ws.Range("F13").NumberFormat = "0"

# Discard that stray bottom row entirely.
$ws.Rows.Item(13).Delete()

# Drop the MARCA / MODELO / TIPO O SERVICIO / CLASE DE SERVICIO columns
# (B:E), shifting CANTIDAD DE PASAJEROS (old column F, with its bestFit
# width) left into column B, and everything after it left as well.
$ws.Range("B1:E1").Delete(-4159)

# Column B (now CANTIDAD DE PASAJEROS) is numeric -> apply the integer
# number format (reuses the style mutated above) while keeping the header
# cell itself on the default look.
$ws.Columns.Item(2).NumberFormat = "0"
$ws.Range("B1").Style = "Normal"

# Give the now-empty column C the width the cleaned-up template uses.
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Move the selection to A2, matching the saved workbook state.
$ws.Range("A2").Select()
